# Updates the cryptos list table (columns B-E, rows 2-51) on Sheet1
# to match the latest scrape: refreshed prices / 1h volume deltas,
# plus the USDe/Cosmos row swap (rows 46-47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.919.08"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "3.928.37"
$ws.Range("E3").Value = "  +4.30%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.37"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.65"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "3.925.59"
$ws.Range("E7").Value = "  +4.27%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.27"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "4.589.13"
$ws.Range("E15").Value = "  +4.43%  "
$ws.Range("D16").Value = "3.934.44"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("D17").Value = "69.047.46"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "487.91"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000166"
$ws.Range("E24").Value = "  +11.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.50"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "4.081.83"
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "3.877.63"
$ws.Range("E35").Value = "  +4.74%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.93"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.323"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "444.72"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.48"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.65"
$ws.Range("E48").Value = "  +13.30%  "
$ws.Range("D49").Value = "2.851.26"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.90"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0359"
$ws.Range("E51").Value = "  +1.74%  "
